$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: advance the header date by one day (17 Jan 2024 -> 18 Jan 2024,
# serial 45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the two unit price cells
$ws.Range("D29").Value = 364.992
$ws.Range("D30").Value = 514.29
